$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for each data row (rows 2-29).
# The diff shows this value being bumped from 45203 (2023-10-04) to 45204 (2023-10-05)
# for every one of those rows, leaving everything else (formatting, other columns) untouched.
for ($row = 2; $row -le 29; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45204
    }
}
